$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2, B2: plain text replacements
$ws.Range("A2").Value = "dsfdsf"
$ws.Range("B2").Value = "vbfvb"

# C2: phone number with a leading "+" — force text storage so Excel doesn't
# silently coerce it into a numeric value and drop the plus sign, then reset
# the cell style back to Normal so no stray number-format style sticks around.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "+380664606978"
$ws.Range("C2").Style = "Normal"

# D2: timestamp string
$ws.Range("D2").Value = "11/12/2023 18:36:45"
